# Auto-generated Excel COM-interop script
# Applies scheduled market-price / profit refresh values to the 8 crafting-leve sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit's data diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 266.1875
$ws.Range("I11").Value = 266.1875
$ws.Range("K11").Value = 266.1875
$ws.Range("M11").Value = -126.1875

$ws.Range("H17").Value = 2505702.5
$ws.Range("J17").Value = 2505702.5
$ws.Range("L17").Value = 7517107.5
$ws.Range("N17").Value = -7517443.5

$ws.Range("H33").Value = 17861122
$ws.Range("I33").Value = 50001300
$ws.Range("K33").Value = 50001300
$ws.Range("M33").Value = -50001071

$ws.Range("H54").Value = 11666.333
$ws.Range("I54").Value = 11666.333
$ws.Range("K54").Value = 11666.333
$ws.Range("M54").Value = -11180.333

$ws.Range("H86").Value = 5310.5264
$ws.Range("I86").Value = 5855.1
$ws.Range("K86").Value = 5855.1
$ws.Range("M86").Value = -4732.1

$ws.Range("H89").Value = 5310.5264
$ws.Range("I89").Value = 5855.1
$ws.Range("K89").Value = 29275.5
$ws.Range("M89").Value = -23659.5

$ws.Range("H98").Value = 2585.0444
$ws.Range("I98").Value = 2277.244
$ws.Range("K98").Value = 2277.244
$ws.Range("M98").Value = -779.2440000000001

$ws.Range("H122").Value = 2585.0444
$ws.Range("I122").Value = 2277.244
$ws.Range("K122").Value = 6831.732
$ws.Range("M122").Value = -4381.732

$ws.Range("H129").Value = 1750.909
$ws.Range("I129").Value = 947.7143
$ws.Range("K129").Value = 2843.1429
$ws.Range("M129").Value = 2156.8571

$ws.Range("H132").Value = 11602.536
$ws.Range("I132").Value = 8873.625
$ws.Range("J132").Value = 27976
$ws.Range("K132").Value = 26620.875
$ws.Range("L132").Value = 83928
$ws.Range("M132").Value = -24090.875
$ws.Range("N132").Value = -88988

$ws.Range("H137").Value = 19022.54
$ws.Range("I137").Value = 1755.2858
$ws.Range("K137").Value = 5265.857400000001
$ws.Range("M137").Value = -2715.857400000001

$ws.Range("H138").Value = 2347.989
$ws.Range("J138").Value = 2479.65
$ws.Range("L138").Value = 7438.950000000001
$ws.Range("N138").Value = -17718.95

$ws.Range("H141").Value = 3821
$ws.Range("I141").Value = 3821
$ws.Range("K141").Value = 11463
$ws.Range("M141").Value = -6283

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4240126
$ws.Range("I32").Value = 4631338.5
$ws.Range("K32").Value = 4631338.5
$ws.Range("M32").Value = -4631051.5

$ws.Range("H61").Value = 55778.883
$ws.Range("I61").Value = 1250
$ws.Range("K61").Value = 1250
$ws.Range("M61").Value = -1038

$ws.Range("H74").Value = 10943.047
$ws.Range("I74").Value = 1431.2222
$ws.Range("K74").Value = 1431.2222
$ws.Range("M74").Value = -557.2221999999999

$ws.Range("H77").Value = 10943.047
$ws.Range("I77").Value = 1431.2222
$ws.Range("K77").Value = 7156.111
$ws.Range("M77").Value = -2788.111

$ws.Range("H102").Value = 3038.5
$ws.Range("I102").Value = 2758.0908
$ws.Range("K102").Value = 2758.0908
$ws.Range("M102").Value = -1136.0908

$ws.Range("H110").Value = 6998606.5
$ws.Range("I110").Value = 7581674
$ws.Range("K110").Value = 7581674
$ws.Range("M110").Value = -7579629

$ws.Range("H132").Value = 2335260.8
$ws.Range("I132").Value = 2404.9143
$ws.Range("J132").Value = 12541505
$ws.Range("K132").Value = 7214.742899999999
$ws.Range("L132").Value = 37624515
$ws.Range("M132").Value = -4684.742899999999
$ws.Range("N132").Value = -37629575

$ws.Range("H136").Value = 55778.883
$ws.Range("I136").Value = 1250
$ws.Range("K136").Value = 3750
$ws.Range("M136").Value = -1200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 9999.5
$ws.Range("J7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("N7").Value = -5226

$ws.Range("H20").Value = 11515686
$ws.Range("I20").Value = 25649680
$ws.Range("K20").Value = 25649680
$ws.Range("M20").Value = -25649433

$ws.Range("H134").Value = 28666.326
$ws.Range("I134").Value = 30877.244
$ws.Range("K134").Value = 92631.73199999999
$ws.Range("M134").Value = -90096.73199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16490.45
$ws.Range("I31").Value = 1137.909
$ws.Range("J31").Value = 35254.668
$ws.Range("K31").Value = 1137.909
$ws.Range("L31").Value = 35254.668
$ws.Range("M31").Value = -842.9090000000001
$ws.Range("N31").Value = -35844.668

$ws.Range("H34").Value = 16490.45
$ws.Range("I34").Value = 1137.909
$ws.Range("J34").Value = 35254.668
$ws.Range("K34").Value = 1137.909
$ws.Range("L34").Value = 35254.668
$ws.Range("M34").Value = -935.9090000000001
$ws.Range("N34").Value = -35658.668

$ws.Range("H105").Value = 11590.5
$ws.Range("I105").Value = 15700.714
$ws.Range("K105").Value = 15700.714
$ws.Range("M105").Value = -13953.714

$ws.Range("H122").Value = 2898.5454
$ws.Range("I122").Value = 1764
$ws.Range("J122").Value = 4260
$ws.Range("K122").Value = 5292
$ws.Range("L122").Value = 12780
$ws.Range("M122").Value = -2842
$ws.Range("N122").Value = -17680

$ws.Range("H125").Value = 85002.336
$ws.Range("J125").Value = 85002.336
$ws.Range("L125").Value = 85002.336
$ws.Range("N125").Value = -89922.336

$ws.Range("H132").Value = 37040370
$ws.Range("J132").Value = 333337150
$ws.Range("L132").Value = 1000011450
$ws.Range("N132").Value = -1000016510

$ws.Range("H134").Value = 55565936
$ws.Range("I134").Value = 2517.625
$ws.Range("J134").Value = 100016670
$ws.Range("K134").Value = 7552.875
$ws.Range("L134").Value = 300050010
$ws.Range("M134").Value = -5017.875
$ws.Range("N134").Value = -300055080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 93.96429000000001
$ws.Range("J2").Value = 40.785713
$ws.Range("L2").Value = 244.714278
$ws.Range("N2").Value = -470.714278

$ws.Range("H8").Value = 17857212
$ws.Range("I8").Value = 17857212
$ws.Range("K8").Value = 53571636
$ws.Range("M8").Value = -53571497

$ws.Range("H98").Value = 2126.4546
$ws.Range("I98").Value = 362
$ws.Range("J98").Value = 3134.7144
$ws.Range("K98").Value = 1086
$ws.Range("L98").Value = 9404.143199999999
$ws.Range("M98").Value = 412
$ws.Range("N98").Value = -12400.1432

$ws.Range("H120").Value = 9000
$ws.Range("I120").Value = 9000
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 27000
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -22162
$ws.Range("N120").ClearContents()

$ws.Range("H131").Value = 1405.49
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 1409.5858
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 4228.7574
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -14308.7574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5660584
$ws.Range("I122").Value = 5660584
$ws.Range("K122").Value = 16981752
$ws.Range("M122").Value = -16979302

$ws.Range("H126").Value = 8954004
$ws.Range("J126").Value = 14289293
$ws.Range("L126").Value = 42867879
$ws.Range("N126").Value = -42872819

$ws.Range("H132").Value = 1507.75
$ws.Range("I132").Value = 1262.125
$ws.Range("K132").Value = 3786.375
$ws.Range("M132").Value = -1256.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1771155.9
$ws.Range("I7").Value = 2983257
$ws.Range("K7").Value = 2983257
$ws.Range("M7").Value = -2983145

$ws.Range("H16").Value = 55556788
$ws.Range("I16").Value = 58824788
$ws.Range("K16").Value = 58824788
$ws.Range("M16").Value = -58824618

$ws.Range("H22").Value = 40001890
$ws.Range("J22").Value = 83335530
$ws.Range("L22").Value = 83335530
$ws.Range("N22").Value = -83336120

$ws.Range("H27").Value = 40001890
$ws.Range("J27").Value = 83335530
$ws.Range("L27").Value = 83335530
$ws.Range("N27").Value = -83335744

$ws.Range("H46").Value = 662.5
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 725
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 725
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -1101

$ws.Range("H126").Value = 1771155.9
$ws.Range("I126").Value = 2983257
$ws.Range("K126").Value = 8949771
$ws.Range("M126").Value = -8947301

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 23405.4
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 26756.75
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 26756.75
$ws.Range("M55").Value = -9723
$ws.Range("N55").Value = -27310.75

$ws.Range("H59").Value = 39967.332
$ws.Range("J59").Value = 39967.332
$ws.Range("L59").Value = 39967.332
$ws.Range("N59").Value = -41443.332

$ws.Range("H81").Value = 3121.3333
$ws.Range("I81").Value = 3501.5386
$ws.Range("K81").Value = 7003.0772
$ws.Range("M81").Value = -5942.0772

$ws.Range("H84").Value = 3121.3333
$ws.Range("I84").Value = 3501.5386
$ws.Range("K84").Value = 35015.386
$ws.Range("M84").Value = -29711.386

$ws.Range("H96").Value = 1887.8182
$ws.Range("I96").Value = 1884.2858
$ws.Range("J96").Value = 1889.4667
$ws.Range("K96").Value = 1884.2858
$ws.Range("L96").Value = 1889.4667
$ws.Range("M96").Value = -511.2858000000001
$ws.Range("N96").Value = -4635.4667
